# Sorted the packet class structure headers to make them more organised.
# Also added the client replying to the server with its start and end positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colorDone = 5287936   # green fill used for "DONE" status cells
$colorTodo = 255       # red fill used for "TODO" status cells

# Row 9: note explaining that the current maze size (as well as the
# user-requested size) should be displayed.
$ws.Range("L9").Value = "Display current maze size as well as user input maze size to be requested"

# Row 10: "Solve memory issue" is now DONE.
$ws.Range("M10").Value = "DONE"
$ws.Range("M10").Interior.Color = $colorDone

# Row 11: client transmits its start and end position to the server.
$ws.Range("K11").Value = "Transmit start and end position"
$ws.Range("L11").Value = "from client to server"
$ws.Range("M11").Value = "DONE"
$ws.Range("M11").Interior.Color = $colorDone

# Row 12: server accepts the position and computes a route.
$ws.Range("K12").Value = "Accept position & compute route"
$ws.Range("L12").Value = "Add options client side for different types of algorithm"
$ws.Range("M12").Value = "TODO"
$ws.Range("M12").Interior.Color = $colorTodo

# Row 13: server sends the route back, client shows it on key press.
$ws.Range("K13").Value = "Send route back to client"
$ws.Range("L13").Value = "Client should display route on key press"
$ws.Range("M13").Value = "TODO"
$ws.Range("M13").Interior.Color = $colorTodo

# Update the sheet's current selection to reflect the cells just edited.
$ws.Range("M10:M11").Select()
